$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# JobOffers sheet: add Level + Skills columns (E:H) for each job offer row
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("JobOffers")

# Header row (row 1) - bold header cells to match existing header styling
$ws2.Range("E1").Value = "Level"
$ws2.Range("F1").Value = "Skills"
$ws2.Range("G1").Value = "Skills"
$ws2.Range("H1").Value = "Skills"
$ws2.Range("E1:H1").Font.Bold = $true

# Row 2
$ws2.Range("E2").Value = "Junior"
$ws2.Range("F2").Value = "Java fundamentals"
$ws2.Range("G2").Value = "Java Spring"

# Row 3
$ws2.Range("E3").Value = "Junior"
$ws2.Range("F3").Value = "Graphics"
$ws2.Range("G3").Value = "Multimedia"

# Row 4
$ws2.Range("E4").Value = "Mid"
$ws2.Range("F4").Value = "Hardware"
$ws2.Range("G4").Value = "Devops"
$ws2.Range("H4").Value = "Databases"

# Row 5
$ws2.Range("E5").Value = "Senior"
$ws2.Range("F5").Value = "Java Spring"

# Row 6
$ws2.Range("E6").Value = "Senior"
$ws2.Range("F6").Value = "Scrum"

# Row 7
$ws2.Range("E7").Value = "Senior"
$ws2.Range("F7").Value = "Problem-solving"

# Row 8
$ws2.Range("E8").Value = "Junior"
$ws2.Range("F8").Value = "Java Spring"
$ws2.Range("G8").Value = "Devops"

# Row 9
$ws2.Range("E9").Value = "Mid"
$ws2.Range("F9").Value = "Operating Systems"
$ws2.Range("G9").Value = "Angular"
$ws2.Range("H9").Value = "C#"

# Row 10
$ws2.Range("E10").Value = "Mid"
$ws2.Range("F10").Value = "Maven"
$ws2.Range("G10").Value = "Spreadsheets"

# Row 11
$ws2.Range("E11").Value = "Junior"
$ws2.Range("F11").Value = "Maven"
$ws2.Range("G11").Value = "Problem-solving"

# Row 12
$ws2.Range("E12").Value = "Junior"
$ws2.Range("F12").Value = "Angular"
$ws2.Range("G12").Value = "C#"

# Row 13
$ws2.Range("E13").Value = "Senior"
$ws2.Range("F13").Value = "Operating Systems"

# Row 14
$ws2.Range("E14").Value = "Mid"
$ws2.Range("F14").Value = "Devops"
$ws2.Range("G14").Value = "Problem-solving"
$ws2.Range("H14").Value = "Databases"

# ---------------------------------------------------------------------------
# Update sheet selections / active views.
# Order matters: the sheet selected/activated last becomes the workbook's
# active tab (tabSelected on its sheetView, activeTab on the workbookView).
# ---------------------------------------------------------------------------

# Applicants sheet: selection becomes G2:G10, no special scroll position
$ws1 = $wb.Worksheets.Item("Applicants")
$ws1.Select()
$ws1.Range("G2:G10").Select()

# Skills sheet: selection becomes A21, no special scroll position, no longer active tab
$ws3 = $wb.Worksheets.Item("Skills")
$ws3.Select()
$ws3.Range("A21").Select()

# JobOffers sheet becomes the active tab with selection E7
$ws2.Select()
$ws2.Range("E7").Select()
